$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.89

$ws.Range("N3").Value = 5.6
$ws.Range("P3").Value = 2.58
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.63
$ws.Range("S3").Value = 2.48
$ws.Range("T3").Value = 1.6
$ws.Range("U3").Value = 2.56
$ws.Range("Y3").Value = 23
$ws.Range("AA3").Value = 95
$ws.Range("AC3").Value = 9.6

$ws.Range("G4").Value = 12.5
$ws.Range("P4").Value = 2.32
$ws.Range("U4").Value = 1.92

$ws.Range("N5").Value = 2.3
$ws.Range("R5").Value = 1.27
$ws.Range("S5").Value = 2.08
$ws.Range("T5").Value = 1.33
$ws.Range("AN5").Value = 13

$ws.Range("F6").Value = 1.39
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 9.4
$ws.Range("J6").Value = 5.8
$ws.Range("K6").Value = 5.9
$ws.Range("O6").Value = 1.16
$ws.Range("R6").Value = 1.78
$ws.Range("S6").Value = 2.2
$ws.Range("T6").Value = 1.77
$ws.Range("V6").Value = 1.12
$ws.Range("W6").Value = 3.5
$ws.Range("Y6").Value = 40
$ws.Range("Z6").Value = 85
$ws.Range("AE6").Value = 110
$ws.Range("AF6").Value = 10.5
$ws.Range("AM6").Value = 95
$ws.Range("AO6").Value = 95

$ws.Range("F7").Value = 3.65
$ws.Range("M7").Value = 1.05
$ws.Range("S7").Value = 2.72
$ws.Range("T7").Value = 1.63
$ws.Range("AB7").Value = 18
$ws.Range("AN7").Value = 29

$ws.Range("F8").Value = 1.69
$ws.Range("G8").Value = 1.7
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 6.4
$ws.Range("J8").Value = 3.95
$ws.Range("N8").Value = 3.75
$ws.Range("S8").Value = 3.6
$ws.Range("T8").Value = 2
$ws.Range("Y8").Value = 19
$ws.Range("AA8").Value = 170
$ws.Range("AF8").Value = 9.2
$ws.Range("AH8").Value = 22
$ws.Range("AJ8").Value = 16
$ws.Range("AN8").Value = 11

$ws.Range("N9").Value = 1.01

$ws.Range("N10").Value = 1.01

$ws.Range("N11").Value = 1.01
$ws.Range("S11").Value = 1.4

$ws.Range("F12").Value = 1.92
$ws.Range("G12").Value = 2.1
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 4.6
$ws.Range("J12").Value = 3.45
$ws.Range("K12").Value = 3.95
$ws.Range("P12").Value = 1.94
$ws.Range("Q12").Value = 1.89
